# configuratie.xlsx - add "verberg lege kolommen" options to the opmaak sheet
#
# Summary of the change (see commit message):
#   "Optie toegevoegd om enkel lege totaalkolommen of lege crossingkolommen
#    te verbergen. Fout opgelost waarbij een subset met een enkele waarde
#    (bijv. regio = GGD) zorgde voor het verdwijnen van de namen."
#
# Concretely, on the "opmaak" worksheet two new configuration rows are
# appended (both booleans, defaulting to FALSE):
#   A30 = "verberg_lege_kolommen_crossing" / B30 = FALSE
#   A31 = "verberg_lege_kolommen_totaal"   / B31 = FALSE
#
# Besides that, the author's last cursor/selection position moved: the
# "datasets" sheet is no longer the active tab (selection now sits on C5
# instead of F2), while "opmaak" becomes the active tab with the selection
# on D28.

$wb = $excel.ActiveWorkbook

# --- "datasets" sheet: just a leftover cursor move, no data changes ------
$wsDatasets = $wb.Worksheets.Item("datasets")
$wsDatasets.Range("C5").Select()

# --- "opmaak" sheet: new boolean settings below the existing ones --------
$wsOpmaak = $wb.Worksheets.Item("opmaak")

$wsOpmaak.Range("A30").Value = "verberg_lege_kolommen_crossing"
$wsOpmaak.Range("B30").Value = $false

$wsOpmaak.Range("A31").Value = "verberg_lege_kolommen_totaal"
$wsOpmaak.Range("B31").Value = $false

# Make "opmaak" the active sheet/tab again, with the cursor parked on D28,
# matching the saved view state from the diff. This must run last so that
# "opmaak" (not "datasets") ends up as the workbook's active tab.
$wsOpmaak.Activate()
$wsOpmaak.Range("D28").Select()
